$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows whose Target cluster is "Neutrophils" (original rows 9 then 5,
# deleted bottom-up so row indices of the remaining rows stay valid while deleting).
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(5).Delete()

# Final data for rows 2-7 (header stays on row 1). Includes the
# "Inflammatory-Mac" -> "Resolving-Mac" cluster rename and refreshed TPM-derived
# metrics for every remaining Sending/Target cluster combination.
$data = @(
    @("FAPs", "Fgf22", "Fgfr2", "ECs", 1, 0.3333333333333333, 0.03143233333333333, 0.094297, 0.3046368955123587, 0.3046368955123587, 1, 0.5, 0.092077, 0.184154, 0.0789959771480734, 0.05545240531440215, 0.002894194956333333, 0.017365169738, 0.02406508923635431, 0.01689284860367249),
    @("FAPs", "Fgf22", "Fgfr2", "FAPs", 1, 0.3333333333333333, 0.03143233333333333, 0.094297, 0.3046368955123587, 0.3046368955123587, 3, 1, 0.9897559999999999, 2.969268, 0.8491451975864605, 0.8941052196698643, 0.03111034051066666, 0.279993064596, 0.2586809568319677, 0.272377438381623),
    @("FAPs", "Fgf22", "Fgfr2", "MuSCs", 1, 0.3333333333333333, 0.03143233333333333, 0.094297, 0.3046368955123587, 0.3046368955123587, 2, 1, 0.083758, 0.167516, 0.07185882526546619, 0.05044237501573352, 0.002632709375333333, 0.015796256252, 0.02189084944403667, 0.01536660852706323),
    @("Resolving-Mac", "Fgf22", "Fgfr2", "ECs", 1, 0.3333333333333333, 0.07174733333333333, 0.215242, 0.6953631044876413, 0.6953631044876413, 1, 0.5, 0.092077, 0.184154, 0.0789959771480734, 0.05545240531440215, 0.006606279211333333, 0.039637675268, 0.05493088791171908, 0.03855955671072965),
    @("Resolving-Mac", "Fgf22", "Fgfr2", "FAPs", 1, 0.3333333333333333, 0.07174733333333333, 0.215242, 0.6953631044876413, 0.6953631044876413, 3, 1, 0.9897559999999999, 2.969268, 0.8491451975864605, 0.8941052196698643, 0.07101235365066666, 0.6391111828559999, 0.5904642407544927, 0.6217277812882414),
    @("Resolving-Mac", "Fgf22", "Fgfr2", "MuSCs", 1, 0.3333333333333333, 0.07174733333333333, 0.215242, 0.6953631044876413, 0.6953631044876413, 2, 1, 0.083758, 0.167516, 0.07185882526546619, 0.05044237501573352, 0.006009413145333333, 0.036056478872, 0.04996797582142953, 0.03507576648867029)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $r).Value = $row[$j]
    }
}
